# wp refactor of tests
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet Q1_20_21
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Q1_20_21")

$ws1.Range("C3").Value = "AMIS"
$ws1.Range("L3").ClearContents()

$ws1.Range("J4").Value = 928
$ws1.Range("L4").Value = "Python is an interpreted, high-level, general-purpose programming language. Created by Guido van Rossum and first released in 1991, Python's design philosophy emphasizes code readability with its notable use of significant whitespace. Its language constructs and object-oriented approach aim to help programmers write clear, logical code for small and large-scale projects.[28] "

$ws1.Range("C5").Value = "RSS"
$ws1.Range("L5").Value = "PyCharm is an integrated development environment (IDE) used in computer programming, specifically for the Python language. It is developed by the Czech company JetBrains.[6] It provides code analysis, a graphical debugger, an integrated unit tester, integration with version control systems (VCSes), and supports web development with Django as well as Data Science with Anaconda.[7] "

$ws1.Range("C6").Value = "RPE"
$ws1.Range("L6").Value = " Datamaps allows for data to be collected from multiple users using pre-defined forms, built out of Excel spreadsheets. The form can be as complex or as simple as you like - Datamaps does the hard work of collecting the data into one place for onward processing or analysis, whilst ensuring that the data is valid according to expectations.  Datamaps acknowledges that in the office environment, we use Excel for everything, but using it to collect data requires help - which is where Datamaps excels. "

$ws1.Range("C7").Value = "RIG"
$ws1.Range("L7").Value = "GitHub, Inc. is an American multinational corporation that provides hosting for software development and version control using Git. It offers the distributed version control and source code management (SCM) functionality of Git, plus its own features. It provides access control and several collaboration features such as bug tracking, feature requests, task management, and wikis for every project.[3] Headquartered in California, it has been a subsidiary of Microsoft since 2018.[4] "

# ---------------------------------------------------------------------
# Sheet Q4_19_20 : row 5 (A11/HSMRPG) is removed and rows 6-8 shift up
# to become rows 5-7, plus the same text edits as above are applied.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Q4_19_20")

$ws2.Range("C3").Value = "AMIS"
$ws2.Range("L3").ClearContents()

$ws2.Range("L4").Value = "Python is an interpreted, high-level, general-purpose programming language. Created by Guido van Rossum and first released in 1991, Python's design philosophy emphasizes code readability with its notable use of significant whitespace. Its language constructs and object-oriented approach aim to help programmers write clear, logical code for small and large-scale projects.[28] "

# Delete row 5 (A11 / HSMRPG) entirely, shifting rows 6-8 up to 5-7.
$ws2.Range("A5").EntireRow.Delete()

$ws2.Range("C5").Value = "RSS"
$ws2.Range("L5").Value = "PyCharm is an integrated development environment (IDE) used in computer programming, specifically for the Python language. It is developed by the Czech company JetBrains.[6] It provides code analysis, a graphical debugger, an integrated unit tester, integration with version control systems (VCSes), and supports web development with Django as well as Data Science with Anaconda.[7] "

$ws2.Range("C6").Value = "RPE"
$ws2.Range("L6").Value = " Datamaps allows for data to be collected from multiple users using pre-defined forms, built out of Excel spreadsheets. The form can be as complex or as simple as you like - Datamaps does the hard work of collecting the data into one place for onward processing or analysis, whilst ensuring that the data is valid according to expectations.  Datamaps acknowledges that in the office environment, we use Excel for everything, but using it to collect data requires help - which is where Datamaps excels. "

$ws2.Range("C7").Value = "RIG"
$ws2.Range("L7").Value = "GitHub, Inc. is an American multinational corporation that provides hosting for software development and version control using Git. It offers the distributed version control and source code management (SCM) functionality of Git, plus its own features. It provides access control and several collaboration features such as bug tracking, feature requests, task management, and wikis for every project.[3] Headquartered in California, it has been a subsidiary of Microsoft since 2018.[4] "

# ---------------------------------------------------------------------
# Sheet Count
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Count")
$ws3.Range("C7").Value = 1761
$ws3.Range("C11").Value = 7978
